$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Route4" column (column E) with origin/destination and timetable data
# (values entered in this order to reproduce the author's shared-string ordering)
$ws.Range("E1").Value = "Route4"
$ws.Range("E2").Value = "'WOK,WAT"
$ws.Range("E3").Value = "'WAT,WOK"
$ws.Range("E4").Value = "'1600,1730,1815"
$ws.Range("E5").Value = "'2300,0000"
$ws.Range("E6").Value = "'1800"
$ws.Range("E8").Value = "'1100,1330"
$ws.Range("E9").Value = "'1215"
$ws.Range("E7").Value = "'0630,0700,0945,1100"

# Set the width of the new column E
$ws.Columns.Item(5).ColumnWidth = 23.3

# Update the active selection to E8
$ws.Range("E8").Select()
